$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 850
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 850
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 850
$ws.Range("M32").ClearContents()
$ws.Range("N32").Value = -1502

$ws.Range("H100").Value = 22224010
$ws.Range("I100").Value = 25001892
$ws.Range("K100").Value = 25001892
$ws.Range("M100").Value = -25001351

$ws.Range("H106").Value = 1763.3334
$ws.Range("I106").Value = 1763.3334
$ws.Range("K106").Value = 1763.3334
$ws.Range("M106").Value = -1132.3334

$ws.Range("H121").Value = 2980
$ws.Range("J121").Value = 2980
$ws.Range("L121").Value = 8940
$ws.Range("N121").Value = -12434

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4131.418
$ws.Range("I32").Value = 4061.0408
$ws.Range("J32").Value = 4706.1665
$ws.Range("K32").Value = 4061.0408
$ws.Range("L32").Value = 4706.1665
$ws.Range("M32").Value = -3774.0408
$ws.Range("N32").Value = -5280.1665

$ws.Range("H33").Value = 18750
$ws.Range("J33").Value = 27500
$ws.Range("L33").Value = 27500
$ws.Range("N33").Value = -28158

$ws.Range("H45").Value = 1622.375
$ws.Range("I45").Value = 1736
$ws.Range("J45").Value = 1433
$ws.Range("K45").Value = 1736
$ws.Range("L45").Value = 1433
$ws.Range("M45").Value = -1359
$ws.Range("N45").Value = -2187

$ws.Range("H48").Value = 79800
$ws.Range("J48").Value = 79800
$ws.Range("L48").Value = 79800
$ws.Range("N48").Value = -80568

$ws.Range("H104").Value = 33500
$ws.Range("J104").Value = 33500
$ws.Range("L104").Value = 33500
$ws.Range("N104").Value = -40488

$ws.Range("H118").Value = 28390
$ws.Range("J118").Value = 28390
$ws.Range("L118").Value = 28390
$ws.Range("N118").Value = -31704

$ws.Range("H122").Value = 2868.4167
$ws.Range("I122").Value = 1920.3334
$ws.Range("J122").Value = 3816.5
$ws.Range("K122").Value = 5761.0002
$ws.Range("L122").Value = 11449.5
$ws.Range("M122").Value = -3311.0002
$ws.Range("N122").Value = -16349.5

$ws.Range("H128").Value = 41445.8
$ws.Range("J128").Value = 41445.8
$ws.Range("L128").Value = 41445.8
$ws.Range("N128").Value = -51405.8

$ws.Range("H137").Value = 40201.668
$ws.Range("J137").Value = 40201.668
$ws.Range("L137").Value = 40201.668
$ws.Range("N137").Value = -50401.668

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H31").Value = 911.5
$ws.Range("I31").Value = 911.5
$ws.Range("K31").Value = 911.5
$ws.Range("M31").Value = -659.5

$ws.Range("H92").Value = 69999.5
$ws.Range("J92").Value = 69999.5
$ws.Range("L92").Value = 69999.5
$ws.Range("N92").Value = -74991.5

$ws.Range("H94").Value = 821.2083
$ws.Range("I94").Value = 648.0476
$ws.Range("J94").Value = 2033.3334
$ws.Range("K94").Value = 648.0476
$ws.Range("L94").Value = 2033.3334
$ws.Range("M94").Value = -197.0476
$ws.Range("N94").Value = -2935.3334

$ws.Range("H137").Value = 43500
$ws.Range("J137").Value = 43500
$ws.Range("L137").Value = 43500
$ws.Range("N137").Value = -53700

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H9").Value = 16946.666
$ws.Range("J9").Value = 16946.666
$ws.Range("L9").Value = 16946.666
$ws.Range("N9").Value = -17282.666

$ws.Range("H58").Value = 2528.6711
$ws.Range("I58").Value = 1610.4364
$ws.Range("J58").Value = 5334.3887
$ws.Range("K58").Value = 1610.4364
$ws.Range("L58").Value = 5334.3887
$ws.Range("M58").Value = -1407.4364
$ws.Range("N58").Value = -5740.3887

$ws.Range("H136").Value = 2528.6711
$ws.Range("I136").Value = 1610.4364
$ws.Range("J136").Value = 5334.3887
$ws.Range("K136").Value = 4831.3092
$ws.Range("L136").Value = 16003.1661
$ws.Range("M136").Value = -2281.3092
$ws.Range("N136").Value = -21103.1661

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 0
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("M92").ClearContents()
$ws.Range("N92").ClearContents()

$ws.Range("H107").Value = 167188.5
$ws.Range("J107").Value = 500600.5
$ws.Range("L107").Value = 1501801.5
$ws.Range("N107").Value = -1505641.5

$ws.Range("H136").Value = 3568.0908
$ws.Range("I136").Value = 3494.4443
$ws.Range("J136").Value = 3899.5
$ws.Range("K136").Value = 10483.3329
$ws.Range("L136").Value = 11698.5
$ws.Range("M136").Value = -5383.332900000001
$ws.Range("N136").Value = -21898.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 35151.2
$ws.Range("J46").Value = 35151.2
$ws.Range("L46").Value = 35151.2
$ws.Range("N46").Value = -35463.2

$ws.Range("H102").Value = 1957.6129
$ws.Range("I102").Value = 1391.8077
$ws.Range("K102").Value = 1391.8077
$ws.Range("M102").Value = 230.1922999999999

$ws.Range("H116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").ClearContents()

$ws.Range("H122").Value = 12202
$ws.Range("I122").Value = 1800
$ws.Range("J122").Value = 15669.333
$ws.Range("K122").Value = 5400
$ws.Range("L122").Value = 47007.999
$ws.Range("M122").Value = -2950
$ws.Range("N122").Value = -51907.999

$ws.Range("H133").Value = 41393.332
$ws.Range("J133").Value = 41393.332
$ws.Range("L133").Value = 41393.332
$ws.Range("N133").Value = -51513.332

$ws.Range("H137").Value = 37222
$ws.Range("J137").Value = 37222
$ws.Range("L137").Value = 37222
$ws.Range("N137").Value = -47422

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H43").Value = 19450
$ws.Range("J43").Value = 19450
$ws.Range("L43").Value = 19450
$ws.Range("N43").Value = -19836

$ws.Range("H114").Value = 38966.668
$ws.Range("J114").Value = 38966.668
$ws.Range("L114").Value = 38966.668
$ws.Range("N114").Value = -47644.668

$ws.Range("H120").Value = 30000
$ws.Range("J120").Value = 30000
$ws.Range("L120").Value = 30000
$ws.Range("N120").Value = -39676

$ws.Range("H122").Value = 6835
$ws.Range("I122").Value = 1000
$ws.Range("J122").Value = 7668.5713
$ws.Range("K122").Value = 3000
$ws.Range("L122").Value = 23005.7139
$ws.Range("M122").Value = -550
$ws.Range("N122").Value = -27905.7139

$ws.Range("H125").Value = 41804.285
$ws.Range("J125").Value = 41804.285
$ws.Range("L125").Value = 41804.285
$ws.Range("N125").Value = -51644.285

$ws.Range("H135").Value = 50000
$ws.Range("J135").Value = 50000
$ws.Range("L135").Value = 50000
$ws.Range("N135").Value = -60140

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H37").Value = 11999.5
$ws.Range("I37").Value = 11999.5
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 11999.5
$ws.Range("L37").Value = 0
$ws.Range("M37").Value = -11796.5
$ws.Range("N37").ClearContents()

$ws.Range("H82").Value = 38307.617
$ws.Range("J82").Value = 38307.617
$ws.Range("L82").Value = 38307.617
$ws.Range("N82").Value = -39073.617

$ws.Range("H85").Value = 38307.617
$ws.Range("J85").Value = 38307.617
$ws.Range("L85").Value = 38307.617
$ws.Range("N85").Value = -40959.617

$ws.Range("H96").Value = 166752320
$ws.Range("I96").Value = 200100780
$ws.Range("K96").Value = 200100780
$ws.Range("M96").Value = -200099407

$ws.Range("H111").Value = 27254.666
$ws.Range("J111").Value = 27254.666
$ws.Range("L111").Value = 27254.666
$ws.Range("N111").Value = -35434.666

$ws.Range("H113").Value = 572.2857
$ws.Range("I113").Value = 572.2857
$ws.Range("K113").Value = 1716.8571
$ws.Range("M113").Value = 453.1428999999998

$ws.Range("H122").Value = 8560
$ws.Range("I122").Value = 3900
$ws.Range("J122").Value = 11666.667
$ws.Range("K122").Value = 11700
$ws.Range("L122").Value = 35000.001
$ws.Range("M122").Value = -9250
$ws.Range("N122").Value = -39900.001

$ws.Range("H123").Value = 38390
$ws.Range("J123").Value = 38390
$ws.Range("L123").Value = 38390
$ws.Range("N123").Value = -48190

$ws.Range("H131").Value = 41812.145
$ws.Range("J131").Value = 41812.145
$ws.Range("L131").Value = 41812.145
$ws.Range("N131").Value = -51892.145

$ws.Range("H136").Value = 4347.96
$ws.Range("I136").Value = 2367.4285
$ws.Range("J136").Value = 6868.636
$ws.Range("K136").Value = 7102.2855
$ws.Range("L136").Value = 20605.908
$ws.Range("M136").Value = -4552.2855
$ws.Range("N136").Value = -25705.908
